# Update Price (D) and Volume 1h (E) columns for the cryptos table
# Numeric-looking Price strings are written with a leading apostrophe so
# Excel stores them as text (matching the original inline-string format)
# instead of silently coercing them into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.714.70'
$ws.Range("E2").Value = '  +2.24%  '
$ws.Range("D3").Value = '1.872.17'
$ws.Range("E3").Value = '  +2.03%  '
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").Value = "'324.44"
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").Value = "'0.4625"
$ws.Range("E7").Value = '  -0.40%  '
$ws.Range("D8").Value = "'0.3861"
$ws.Range("E8").Value = '  -0.21%  '
$ws.Range("D9").Value = "'0.07874"
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("D10").Value = "'0.9761"
$ws.Range("E10").Value = '  +1.56%  '
$ws.Range("D11").Value = "'21.83"
$ws.Range("E11").Value = '  -0.53%  '
$ws.Range("D12").Value = '1.866.10'
$ws.Range("E12").Value = '  +3.36%  '
$ws.Range("D13").Value = "'7.021"
$ws.Range("E13").Value = '  +1.61%  '
$ws.Range("E14").Value = '  +0.22%  '
$ws.Range("D15").Value = "'0.06952"
$ws.Range("E15").Value = '  +1.43%  '
$ws.Range("D16").Value = "'88.44"
$ws.Range("E16").Value = '  +1.44%  '
$ws.Range("D17").Value = "'1.004"
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("D18").Value = "'0.00001003"
$ws.Range("E18").Value = '  +0.90%  '
$ws.Range("D19").Value = "'16.81"
$ws.Range("E19").Value = '  +0.88%  '
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").Value = '28.689.69'
$ws.Range("E21").Value = '  +2.12%  '
$ws.Range("D22").Value = "'5.275"
$ws.Range("E22").Value = '  -0.99%  '
$ws.Range("D23").Value = "'11.09"
$ws.Range("E23").Value = '  +0.84%  '
$ws.Range("D24").Value = "'2.101"
$ws.Range("E24").Value = '  +0.25%  '
$ws.Range("D25").Value = '2.069.60'
$ws.Range("E25").Value = '  +1.76%  '
$ws.Range("D26").Value = "'152.96"
$ws.Range("E26").Value = '  -0.74%  '
$ws.Range("D27").Value = "'19.30"
$ws.Range("E27").Value = '  +0.64%  '
$ws.Range("D28").Value = "'5.886"
$ws.Range("E28").Value = '  +3.35%  '
$ws.Range("D29").Value = "'1.985"
$ws.Range("E29").Value = '  +0.91%  '
$ws.Range("D30").Value = "'119.27"
$ws.Range("E30").Value = '  +1.17%  '
$ws.Range("D31").Value = "'0.09332"
$ws.Range("E31").Value = '  +0.84%  '
$ws.Range("D32").Value = "'0.9187"
$ws.Range("E32").Value = '  -1.89%  '
$ws.Range("D33").Value = "'5.289"
$ws.Range("E33").Value = '  +0.30%  '
$ws.Range("D34").Value = "'1.335"
$ws.Range("E34").Value = '  +0.87%  '
$ws.Range("D35").Value = "'3.324"
$ws.Range("E35").Value = '  +0.96%  '
$ws.Range("D36").Value = "'0.05798"
$ws.Range("E36").Value = '  -0.90%  '
$ws.Range("D37").Value = "'1.155"
$ws.Range("E37").Value = '  +1.02%  '
$ws.Range("D38").Value = "'0.02077"
$ws.Range("E38").Value = '  -2.26%  '
$ws.Range("D39").Value = "'7.665"
$ws.Range("E39").Value = '  -1.47%  '
$ws.Range("D40").Value = "'0.5623"
$ws.Range("E40").Value = '  +0.41%  '
$ws.Range("D41").Value = "'0.1783"
$ws.Range("E41").Value = '  +1.04%  '
$ws.Range("D42").Value = "'9.775"
$ws.Range("E42").Value = '  -1.24%  '
$ws.Range("D43").Value = "'0.07220"
$ws.Range("E43").Value = '  +0.06%  '
$ws.Range("D44").Value = "'11.77"
$ws.Range("E44").Value = '  +1.18%  '
$ws.Range("D45").Value = "'0.5286"
$ws.Range("E45").Value = '  +0.27%  '
$ws.Range("D46").Value = "'2.146"
$ws.Range("E46").Value = '  +1.24%  '
$ws.Range("D47").Value = "'1.117"
$ws.Range("E47").Value = '  -0.33%  '
$ws.Range("D48").Value = "'1.837"
$ws.Range("E48").Value = '  +0.10%  '
$ws.Range("D49").Value = "'112.84"
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").Value = "'2.416"
$ws.Range("E50").Value = '  +4.12%  '
$ws.Range("D51").Value = "'1.002"
